$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table (with header in row 1) lives in A1:F11. Sort the body
# (A2:F11) first by modelType (column A), then by transType (column B),
# both ascending - mirroring a manual Data > Sort the author did before
# writing the parameter-estimate table out to file.
$dataRange = $ws.Range("A2:F11")
$key1 = $ws.Range("A2:A11")
$key2 = $ws.Range("B2:B11")

$dataRange.Sort($key1, 1, $key2, [Type]::Missing, 1, [Type]::Missing, 1, 0)

# Re-enter the totalFreeRateParams formulas row by row so each cell holds
# its own SUM formula (rather than one shared across the rows) now that
# the rows have been reordered by the sort.
for ($row = 2; $row -le 11; $row++) {
    $ws.Range("F" + $row).Formula = "=SUM(D" + $row + ":E" + $row + ")"
}

# Leave the numeric columns used for the sort selected, matching the
# selection state left behind after performing the sort.
$ws.Range("D2:E11").Select()
